# Insert a new data row at row 139 (a new daily price record), shifting all
# subsequent rows (old 139..219) down by one to (140..220).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(139).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A139").Value = 10
$ws.Range("B139").Value = "Vega Modelo de Temuco"
$ws.Range("C139").Value = "La Araucanía"
$ws.Range("D139").Value = 44438
$ws.Range("E139").Value = 9
$ws.Range("F139").Value = 100112032
$ws.Range("G139").Value = "Zapallo italiano"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 130
$ws.Range("K139").Value = 15000
$ws.Range("L139").Value = 17000
$ws.Range("M139").Value = 16231
$ws.Range("N139").Value = "`$/caja 60 unidades"
$ws.Range("O139").Value = "Región de Arica y Parinacota"
$ws.Range("P139").Value = 271
$ws.Range("Q139").Value = 60
$ws.Range("R139").Value = "Hortaliza"
